$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from bottom to top so row indices remain valid during the operation
$ws.Rows.Item(62).Delete()
$ws.Rows.Item(61).Delete()
$ws.Rows.Item(24).Delete()
